$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename header row labels: "..._old" -> "..._FV2410", "..._new" -> "..._FV2504"
#    (column K stays "diff")
# ---------------------------------------------------------------------------
$fv2410Headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)
$fv2504Headers = @(
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

$leftCols = @("A","B","C","D","E","F","G","H","I","J")
for ($i = 0; $i -lt $leftCols.Length; $i++) {
    $ws.Range($leftCols[$i] + "1").Value = $fv2410Headers[$i]
}

$rightCols = @("L","M","N","O","P","Q","R","S","T","U")
for ($i = 0; $i -lt $rightCols.Length; $i++) {
    $ws.Range($rightCols[$i] + "1").Value = $fv2504Headers[$i]
}

# ---------------------------------------------------------------------------
# 2) Turn the used range A1:U79 into an Excel Table ("Table1"), without
#    letting the table feature capture the existing bold/shaded header
#    formatting as a one-off dxf override, and without leaving a named
#    table style behind.
# ---------------------------------------------------------------------------
$usedRange = $ws.Range("A1:U79")
$headerRange = $ws.Range("A1:U1")

# Stash the header row's current formatting on a scratch row far below the
# data, so it can be restored after the table is created (table creation
# otherwise records the header's pre-existing formatting as a dxf and wires
# it up via headerRowDxfId).
$scratch = $ws.Range("A200:U200")
$headerRange.Copy()
$scratch.PasteSpecial(-4122)
$excel.CutCopyMode = 0

$headerRange.Style = "Normal"

$tbl = $ws.ListObjects.Add(1, $usedRange, [System.Reflection.Missing]::Value, 1)
$tbl.TableStyle = ""

$scratch.Copy()
$headerRange.PasteSpecial(-4122)
$excel.CutCopyMode = 0
$scratch.Clear()

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split after row 1).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()
